$d = $word.ActiveDocument

# The document has three inline pictures living in header/footer stories:
#   - Section 1, Header(wdHeaderFooterFirstPage=2): "BTec_Logo-Orange"  -> rename image1.jpg -> image2.jpg
#   - Section 1, Footer(wdHeaderFooterPrimary=1):   "...PearsonLogo.png" -> rename image2.png -> image1.png
#   - Section 1, Footer(wdHeaderFooterFirstPage=2): "...PearsonLogo.png" -> rename image2.png -> image1.png
#
# Walk every section/header/footer so the script is resilient to section
# count, rather than hard-coding indices blindly.

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $sec = $d.Sections.Item($s)

    for ($i = 1; $i -le $sec.Headers.Count; $i++) {
        $hdr = $sec.Headers.Item($i)
        if ($hdr.Exists) {
            $rng = $hdr.Range
            for ($j = 1; $j -le $rng.InlineShapes.Count; $j++) {
                $shp = $rng.InlineShapes.Item($j)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image2.jpg"
                }
            }
        }
    }

    for ($i = 1; $i -le $sec.Footers.Count; $i++) {
        $ftr = $sec.Footers.Item($i)
        if ($ftr.Exists) {
            $rng = $ftr.Range
            for ($j = 1; $j -le $rng.InlineShapes.Count; $j++) {
                $shp = $rng.InlineShapes.Item($j)
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image1.png"
                }
            }
        }
    }
}
